# TC35_Canine_Filter_Breed-Miniature.xlsx - "corrected ICDC Breed 1-14 scripts"
#
# The FilesTab Cypher query (cell B4 on the "startup" sheet) had two
# columns removed from its RETURN clause:
#   - coalesce(f.file_type, '') AS `File Type`
#   - coalesce(demo.breed,'') AS Breed
# which also shortens the wrapped text by two lines, shrinking row 4's
# auto-fit height from 246.5 to 217.5 points. The active selection on the
# sheet was also left on B4 (the cell that was edited) instead of D11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed  IN ['Miniature Schnauzer'] 
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS `File Name`,
         coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

$ws.Range("B4").Value = $newQuery

# Row 4 wraps B4's text, so removing two lines shrinks its auto height.
$ws.Rows.Item(4).RowHeight = 217.5

# Leave the selection on the cell that was edited.
$ws.Range("B4").Select()
